$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '44.100.79'
$ws.Range('E2').Value = '  +3.51%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.246.26'
$ws.Range('E3').Value = '  +2.18%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
Set-TextValue $ws.Range('D5') '257.98'
$ws.Range('E5').Value = '  +2.73%  '

# Row 6
Set-TextValue $ws.Range('D6') '80.33'
$ws.Range('E6').Value = '  +7.76%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.628'
$ws.Range('E7').Value = '  +3.59%  '

# Row 8
$ws.Range('E8').Value = '  +0.09%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.604'
$ws.Range('E9').Value = '  +2.75%  '

# Row 10
Set-TextValue $ws.Range('D10') '43.66'
$ws.Range('E10').Value = '  +8.42%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0938'
$ws.Range('E11').Value = '  +1.86%  '

# Row 12
Set-TextValue $ws.Range('D12') '7.11'
$ws.Range('E12').Value = '  +4.03%  '

# Row 13
$ws.Range('E13').Value = '  +2.06%  '

# Row 14
Set-TextValue $ws.Range('D14') '2.588.78'
$ws.Range('E14').Value = '  +2.31%  '

# Row 15
Set-TextValue $ws.Range('D15') '14.77'
$ws.Range('E15').Value = '  +3.10%  '

# Row 16
Set-TextValue $ws.Range('D16') '2.272.54'
$ws.Range('E16').Value = '  +3.66%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.793'
$ws.Range('E17').Value = '  +1.59%  '

# Row 18
Set-TextValue $ws.Range('D18') '44.058.16'

# Row 19
$ws.Range('E19').Value = '  +2.79%  '

# Row 20
Set-TextValue $ws.Range('D20') '71.66'
$ws.Range('E20').Value = '  +0.70%  '

# Row 21
Set-TextValue $ws.Range('D21') '6.08'
$ws.Range('E21').Value = '  +2.56%  '

# Row 22
Set-TextValue $ws.Range('D22') '2.35'
$ws.Range('E22').Value = '  +9.08%  '

# Row 23
Set-TextValue $ws.Range('D23') '234.66'
$ws.Range('E23').Value = '  +2.68%  '

# Row 24
Set-TextValue $ws.Range('D24') '9.37'
$ws.Range('E24').Value = '  -0.80%  '

# Row 25
$ws.Range('E25').Value = '  +0.15%  '

# Row 26
Set-TextValue $ws.Range('D26') '10.91'
$ws.Range('E26').Value = '  +1.95%  '

# Row 27
Set-TextValue $ws.Range('D27') '40.64'
$ws.Range('E27').Value = '  +8.65%  '

# Row 28
Set-TextValue $ws.Range('D28') '3.37'
$ws.Range('E28').Value = '  -0.04%  '

# Row 29
Set-TextValue $ws.Range('D29') '2.24'
$ws.Range('E29').Value = '  +2.04%  '

# Row 30
$ws.Range('E30').Value = '  -0.68%  '

# Row 31
Set-TextValue $ws.Range('D31') '173.12'

# Row 32
Set-TextValue $ws.Range('D32') '20.65'
$ws.Range('E32').Value = '  +2.82%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.0879'
$ws.Range('E33').Value = '  +10.11%  '

# Row 34
Set-TextValue $ws.Range('D34') '5.33'
$ws.Range('E34').Value = '  +3.05%  '

# Row 35
$ws.Range('E35').Value = '  +5.91%  '

# Row 36
$ws.Range('E36').Value = '  +1.83%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.0368'
$ws.Range('E37').Value = '  +12.50%  '

# Row 38
Set-TextValue $ws.Range('D38') '4.52'
$ws.Range('E38').Value = '  +3.37%  '

# Row 39
Set-TextValue $ws.Range('D39') '12.97'
$ws.Range('E39').Value = '  +6.48%  '

# Row 40
Set-TextValue $ws.Range('D40') '2.89'
$ws.Range('E40').Value = '  +18.94%  '

# Row 41
Set-TextValue $ws.Range('D41') '2.15'
$ws.Range('E41').Value = '  +3.32%  '

# Row 42
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range('D42') '63.36'
$ws.Range('E42').Value = '  +6.94%  '

# Row 43
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D43') '5.54'
$ws.Range('E43').Value = '  +5.15%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.204'
$ws.Range('E44').Value = '  +2.88%  '

# Row 45
Set-TextValue $ws.Range('D45') '104.29'
$ws.Range('E45').Value = '  +1.07%  '

# Row 46
Set-TextValue $ws.Range('D46') '8.57'
$ws.Range('E46').Value = '  +1.16%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.0991'
$ws.Range('E47').Value = '  +1.22%  '

# Row 48
Set-TextValue $ws.Range('D48') '0.458'
$ws.Range('E48').Value = '  -4.78%  '

# Row 49
$ws.Range('E49').Value = '  +2.82%  '

# Row 50
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D50') '1.16'
$ws.Range('E50').Value = '  +2.49%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D51') '1.52'
$ws.Range('E51').Value = '  +23.88%  '
